$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 7313.7144
$ws.Range("J2").Value = 182
$ws.Range("L2").Value = 182
$ws.Range("N2").Value = -408
$ws.Range("H28").Value = 33719.03
$ws.Range("I28").Value = 40918.24
$ws.Range("K28").Value = 40918.24
$ws.Range("M28").Value = -40433.24
$ws.Range("H62").Value = 680391.9399999999
$ws.Range("I62").Value = 1474766.4
$ws.Range("J62").Value = 62545.11
$ws.Range("K62").Value = 1474766.4
$ws.Range("L62").Value = 62545.11
$ws.Range("M62").Value = -1474142.4
$ws.Range("N62").Value = -63793.11
$ws.Range("H65").Value = 680391.9399999999
$ws.Range("I65").Value = 1474766.4
$ws.Range("J65").Value = 62545.11
$ws.Range("K65").Value = 7373832
$ws.Range("L65").Value = 312725.55
$ws.Range("M65").Value = -7370712
$ws.Range("N65").Value = -318965.55
$ws.Range("H98").Value = 2018.75
$ws.Range("I98").Value = 2018.75
$ws.Range("K98").Value = 2018.75
$ws.Range("M98").Value = -520.75
$ws.Range("H111").Value = 27617.334
$ws.Range("I111").Value = 10705
$ws.Range("K111").Value = 32115
$ws.Range("M111").Value = -29048
$ws.Range("H118").Value = 204.1
$ws.Range("I118").Value = 204.1
$ws.Range("K118").Value = 612.3
$ws.Range("M118").Value = 1044.7
$ws.Range("H122").Value = 2018.75
$ws.Range("I122").Value = 2018.75
$ws.Range("K122").Value = 6056.25
$ws.Range("M122").Value = -3606.25
$ws.Range("H132").Value = 3360.1807
$ws.Range("I132").Value = 3223.5166
$ws.Range("K132").Value = 9670.549800000001
$ws.Range("M132").Value = -7140.549800000001
$ws.Range("H135").Value = 62500740
$ws.Range("I135").Value = 66667388
$ws.Range("J135").Value = 1037
$ws.Range("K135").Value = 600006492
$ws.Range("L135").Value = 9333
$ws.Range("M135").Value = -600003957
$ws.Range("N135").Value = -14403
$ws.Range("H136").Value = 60000
$ws.Range("J136").Value = 60000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200
$ws.Range("H138").Value = 3564.8
$ws.Range("J138").Value = 4180.6875
$ws.Range("L138").Value = 12542.0625
$ws.Range("N138").Value = -22822.0625

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 50000000
$ws.Range("J92").Value = 50000000
$ws.Range("L92").Value = 50000000
$ws.Range("N92").Value = -50004992
$ws.Range("H110").Value = 34488800
$ws.Range("I110").Value = 38462356
$ws.Range("K110").Value = 38462356
$ws.Range("M110").Value = -38460311
$ws.Range("H132").Value = 22224532
$ws.Range("I132").Value = 27029198
$ws.Range("K132").Value = 81087594
$ws.Range("M132").Value = -81085064

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1693.4572
$ws.Range("I86").Value = 1613.3
$ws.Range("J86").Value = 1800.3334
$ws.Range("K86").Value = 1613.3
$ws.Range("L86").Value = 1800.3334
$ws.Range("M86").Value = -490.3
$ws.Range("N86").Value = -4046.3334
$ws.Range("H89").Value = 1693.4572
$ws.Range("I89").Value = 1613.3
$ws.Range("J89").Value = 1800.3334
$ws.Range("K89").Value = 8066.5
$ws.Range("L89").Value = 9001.666999999999
$ws.Range("M89").Value = -2450.5
$ws.Range("N89").Value = -20233.667
$ws.Range("H107").Value = 21757362
$ws.Range("I107").Value = 11296.611
$ws.Range("K107").Value = 11296.611
$ws.Range("M107").Value = -9376.611000000001
$ws.Range("H141").Value = 68333.11
$ws.Range("J141").Value = 68333.11
$ws.Range("L141").Value = 68333.11
$ws.Range("N141").Value = -78693.11

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2479.2
$ws.Range("I4").Value = 195
$ws.Range("K4").Value = 195
$ws.Range("M4").Value = -83
$ws.Range("H16").Value = 1575.375
$ws.Range("I16").Value = 1636.2727
$ws.Range("J16").Value = 1441.4
$ws.Range("K16").Value = 1636.2727
$ws.Range("L16").Value = 1441.4
$ws.Range("M16").Value = -1349.2727
$ws.Range("N16").Value = -2015.4
$ws.Range("H31").Value = 4872.857
$ws.Range("I31").Value = 3919.25
$ws.Range("J31").Value = 5097.2354
$ws.Range("K31").Value = 3919.25
$ws.Range("L31").Value = 5097.2354
$ws.Range("M31").Value = -3624.25
$ws.Range("N31").Value = -5687.2354
$ws.Range("H34").Value = 4872.857
$ws.Range("I34").Value = 3919.25
$ws.Range("J34").Value = 5097.2354
$ws.Range("K34").Value = 3919.25
$ws.Range("L34").Value = 5097.2354
$ws.Range("M34").Value = -3717.25
$ws.Range("N34").Value = -5501.2354
$ws.Range("H58").Value = 2229.2104
$ws.Range("I58").Value = 1310.5
$ws.Range("J58").Value = 3250
$ws.Range("K58").Value = 1310.5
$ws.Range("L58").Value = 3250
$ws.Range("M58").Value = -1107.5
$ws.Range("N58").Value = -3656
$ws.Range("H113").Value = 1575.375
$ws.Range("I113").Value = 1636.2727
$ws.Range("J113").Value = 1441.4
$ws.Range("K113").Value = 1636.2727
$ws.Range("L113").Value = 1441.4
$ws.Range("M113").Value = 533.7273
$ws.Range("N113").Value = -5781.4
$ws.Range("H133").Value = 40000
$ws.Range("I133").Value = 40000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 40000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -37470
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 2922.0715
$ws.Range("I134").Value = 1901.6666
$ws.Range("K134").Value = 5704.9998
$ws.Range("M134").Value = -3169.9998
$ws.Range("H136").Value = 2229.2104
$ws.Range("I136").Value = 1310.5
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 3931.5
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -1381.5
$ws.Range("N136").Value = -14850

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 567.7143
$ws.Range("I55").Value = 354.8
$ws.Range("J55").Value = 1100
$ws.Range("K55").Value = 1064.4
$ws.Range("L55").Value = 3300
$ws.Range("M55").Value = -887.4000000000001
$ws.Range("N55").Value = -3654
$ws.Range("H68").Value = 1985.5555
$ws.Range("I68").Value = 1625
$ws.Range("J68").Value = 2165.8333
$ws.Range("K68").Value = 4875
$ws.Range("L68").Value = 6497.499899999999
$ws.Range("M68").Value = -4064
$ws.Range("N68").Value = -8119.499899999999
$ws.Range("H71").Value = 1985.5555
$ws.Range("I71").Value = 1625
$ws.Range("J71").Value = 2165.8333
$ws.Range("K71").Value = 14625
$ws.Range("L71").Value = 19492.4997
$ws.Range("M71").Value = -10569
$ws.Range("N71").Value = -27604.4997
$ws.Range("H106").Value = 3002
$ws.Range("J106").Value = 3002
$ws.Range("L106").Value = 9006
$ws.Range("N106").Value = -10898
$ws.Range("H121").Value = 84903.664
$ws.Range("I121").Value = 872.3333
$ws.Range("K121").Value = 2616.9999
$ws.Range("M121").Value = -1306.9999
$ws.Range("H139").Value = 1239956.6
$ws.Range("I139").Value = 1760827.9
$ws.Range("K139").Value = 5282483.699999999
$ws.Range("M139").Value = -5277343.699999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11731.625
$ws.Range("I70").Value = 12200.889
$ws.Range("J70").Value = 11128.286
$ws.Range("K70").Value = 12200.889
$ws.Range("L70").Value = 11128.286
$ws.Range("M70").Value = -11930.889
$ws.Range("N70").Value = -11668.286
$ws.Range("H73").Value = 11731.625
$ws.Range("I73").Value = 12200.889
$ws.Range("J73").Value = 11128.286
$ws.Range("K73").Value = 12200.889
$ws.Range("L73").Value = 11128.286
$ws.Range("M73").Value = -11264.889
$ws.Range("N73").Value = -13000.286
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H111").Value = 20000
$ws.Range("I111").Value = 20000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 20000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -16933
$ws.Range("N111").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("H122").Value = 36113610
$ws.Range("I122").Value = 465459.25
$ws.Range("J122").Value = 250002500
$ws.Range("K122").Value = 1396377.75
$ws.Range("L122").Value = 750007500
$ws.Range("M122").Value = -1393927.75
$ws.Range("N122").Value = -750012400
$ws.Range("H132").Value = 405003.4
$ws.Range("I132").Value = 629388.5
$ws.Range("J132").Value = 6096.5557
$ws.Range("K132").Value = 1888165.5
$ws.Range("L132").Value = 18289.6671
$ws.Range("M132").Value = -1885635.5
$ws.Range("N132").Value = -23349.6671

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4745581.5
$ws.Range("I2").Value = 19
$ws.Range("J2").Value = 11863925
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 11863925
$ws.Range("M2").Value = 93
$ws.Range("N2").Value = -11864149
$ws.Range("H46").Value = 2115.2415
$ws.Range("I46").Value = 1945.3226
$ws.Range("K46").Value = 1945.3226
$ws.Range("M46").Value = -1757.3226
$ws.Range("H82").Value = 1470.875
$ws.Range("I82").Value = 2946
$ws.Range("J82").Value = 979.1667
$ws.Range("K82").Value = 2946
$ws.Range("L82").Value = 979.1667
$ws.Range("M82").Value = -2585
$ws.Range("N82").Value = -1701.1667
$ws.Range("H85").Value = 1470.875
$ws.Range("I85").Value = 2946
$ws.Range("J85").Value = 979.1667
$ws.Range("K85").Value = 2946
$ws.Range("L85").Value = 979.1667
$ws.Range("M85").Value = -1698
$ws.Range("N85").Value = -3475.1667
$ws.Range("H132").Value = 4645.7837
$ws.Range("I132").Value = 2551
$ws.Range("J132").Value = 7718.1333
$ws.Range("K132").Value = 7653
$ws.Range("L132").Value = 23154.3999
$ws.Range("M132").Value = -5123
$ws.Range("N132").Value = -28214.3999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1812.1482
$ws.Range("I122").Value = 1678.4706
$ws.Range("J122").Value = 2039.4
$ws.Range("K122").Value = 2039.4
$ws.Range("L122").Value = 6118.200000000001
$ws.Range("M122").Value = -2585.4118
$ws.Range("N122").Value = -11018.2
